# Weekly update: a new week's price observation for Brocoli at
# "Macroferia Regional de Talca" is inserted as the new row 550.
# All existing historical rows 550..662 shift down by one row
# (551..663), and the data that used to be on row 662 becomes row 663.
#
# Only columns D, I, J, K, L, M, O, P actually vary row-to-row for this
# data block; A, B, C, E, F, G, H, N, Q, R are constant across the whole
# block, so we simply re-write every column for the rows that move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastOldRow = 662
$insertAt = 550
$newLastRow = $lastOldRow + 1

# Capture the full original last row (A..R) before anything is
# overwritten - it becomes the new row 663.
$savedLastRow = @{}
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $allCols) {
    $savedLastRow[$col] = $ws.Range($col + $lastOldRow).Value()
}

# Columns that actually differ row-to-row in this data block.
$varCols = @("D","I","J","K","L","M","O","P")

# Shift rows down by one: starting from the bottom so we never
# clobber a source row before it has been read.
for ($r = $lastOldRow; $r -ge $insertAt + 1; $r--) {
    $srcRow = $r - 1
    foreach ($col in $varCols) {
        $val = $ws.Range($col + $srcRow).Value()
        $ws.Range($col + $r).Value = $val
    }
}

# The brand-new observation week goes into row 550; everything except
# the date (D) stays the same as what used to be there.
$ws.Range("D" + $insertAt).Value = 45244

# Re-create the constant columns on the new last row and restore the
# values captured at the start (this is what used to be row 662).
foreach ($col in $allCols) {
    $ws.Range($col + $newLastRow).Value = $savedLastRow[$col]
}
